# Apply the edit: in column A (session_name), replace every occurrence of
# ": " (colon followed by a space) with "-" (a single hyphen), for all data
# rows (everything below the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains(": ")) {
        $cell.Value = $val.Replace(": ", "-")
    }
}
